$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(74).Insert()

$ws.Cells.Item(74, 1).Value = 5
$ws.Cells.Item(74, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(74, 3).Value = "Maule"
$ws.Cells.Item(74, 4).Value = 44889
$ws.Cells.Item(74, 5).Value = 7
$ws.Cells.Item(74, 6).Value = 300000000
$ws.Cells.Item(74, 7).Value = "Espárragos"
$ws.Cells.Item(74, 8).Value = "Sin especificar"
$ws.Cells.Item(74, 9).Value = "Primera"
$ws.Cells.Item(74, 10).Value = 3000
$ws.Cells.Item(74, 11).Value = 1000
$ws.Cells.Item(74, 12).Value = 1000
$ws.Cells.Item(74, 13).Value = 1000
$ws.Cells.Item(74, 14).Value = "$/kilo"
$ws.Cells.Item(74, 15).Value = "Provincia de Linares"
$ws.Cells.Item(74, 16).Value = 1000
$ws.Cells.Item(74, 17).Value = 1
$ws.Cells.Item(74, 18).Value = "Hortaliza"
